$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new log entries (rows 29 and 30) ---
# Copy formatting from the last existing data row (28) into the new rows.
$ws.Range("A28:F28").Copy($ws.Range("A29:F29"))
$ws.Range("A28:F28").Copy($ws.Range("A30:F30"))

# Row 29: Constraints chapter work
$ws.Range("A29").Value = "14.10.2022"
$ws.Range("B29").Value = 0.41666666666666669
$ws.Range("C29").Value = "Constraints"
$ws.Range("D29").Value = "Documentation"
$ws.Range("E29").Value = 55
$ws.Range("F29").Value = "Constraints Chapter: Scope and Limitations"

# Row 30: Cost, Marketability, Feasibility chapters work
$ws.Range("A30").Value = "14.10.2022"
$ws.Range("B30").Value = 0.58333333333333337
$ws.Range("C30").Value = "Cost, Marketability, Feasibility"
$ws.Range("D30").Value = "Documentation"
$ws.Range("E30").Value = 220
$ws.Range("F30").Value = "Chapters Cost (research prototype parts for estimation), Target Audience, Marketability, Feasibility and Social Aspects for Success"

# --- Move the summary block (Minutes / Hours / footer) down from rows 31-33 to 39-41 ---
$ws.Range("D31:E32").Copy($ws.Range("D39:E40"))
$ws.Range("C33:F33").Copy($ws.Range("C41:F41"))
$ws.Range("A31:F33").Clear()

# Recompute totals to include the two new rows.
$ws.Range("E39").Formula = "=SUM(E2:E30)"
$ws.Range("E40").Formula = "=E39 / 60"

# Match the saved selection state.
[void]$ws.Range("E40").Select()
